# Adapt column header formatting to respective input file names (#7)
#
# - Renames the "_old" / "_new" header suffixes to "_FV2404" / "_FV2410"
# - Wraps the data range in an Excel Table ("Table1")
# - Freezes the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A-J: "<name>_old" -> "<name>_FV2404"
for ($i = 0; $i -lt $oldNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $oldNames[$i] + "_FV2404"
}

# Column K stays "diff" (column 11) - unchanged

# Columns L-U: "<name>_new" -> "<name>_FV2410"
for ($i = 0; $i -lt $oldNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $oldNames[$i] + "_FV2410"
}

# Turn the used range into an Excel Table named "Table1"
$tableRange = $ws.Range("A1:U65")
$lo = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
$lo.ShowTableStyleRowStripes = $true
$lo.ShowTableStyleColumnStripes = $false
$lo.ShowTableStyleFirstColumn = $false
$lo.ShowTableStyleLastColumn = $false

# Freeze the header row (split below row 1, freeze panes)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

